# Apply the "gh-pages output generated at 456a3b4" update to
# 北京-漫展信息.xlsx: refresh "想去人数" (want-to-go count) figures on the
# 展览 (Exhibitions) and 全部类型 (All types) sheets, plus rename the
# "北京·原神only" event to "北京·原神only3.0" on both sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F4").Value = 3647
$ws1.Range("F5").Value = 3647
$ws1.Range("F6").Value = 273
$ws1.Range("F7").Value = 5175
$ws1.Range("F8").Value = 549
$ws1.Range("F13").Value = 105
$ws1.Range("F15").Value = 713
$ws1.Range("F16").Value = 324
$ws1.Range("F22").Value = 4951
$ws1.Range("F29").Value = 3233
$ws1.Range("F30").Value = 350
$ws1.Range("F31").Value = 720
$ws1.Range("F34").Value = 127
$ws1.Range("F36").Value = 1062
$ws1.Range("F37").Value = 85
$ws1.Range("F41").Value = 1040
$ws1.Range("C42").Value = "北京·原神only3.0"

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F7").Value = 3647
$ws4.Range("F8").Value = 3647
$ws4.Range("F9").Value = 273
$ws4.Range("F10").Value = 5175
$ws4.Range("F11").Value = 549
$ws4.Range("F16").Value = 105
$ws4.Range("F18").Value = 713
$ws4.Range("F19").Value = 324
$ws4.Range("F26").Value = 4951
$ws4.Range("F33").Value = 3233
$ws4.Range("F34").Value = 350
$ws4.Range("F35").Value = 720
$ws4.Range("F39").Value = 127
$ws4.Range("F41").Value = 1062
$ws4.Range("F42").Value = 85
$ws4.Range("F46").Value = 1040
$ws4.Range("C48").Value = "北京·原神only3.0"
